# Remove the "Draft 2" and "Draft 3" sample sections (and everything
# after them, up to the end of the document body) that followed the
# "No hosting plan yet for future backend" bullet in Draft 1. Those
# drafts were unwanted duplicate/experimental content.

$d = $word.ActiveDocument

# Locate the last bullet of the content we want to KEEP (end of Draft 1).
$anchor = $d.Content
$anchor.Find.Execute("No hosting plan yet for future backend", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)

# Resolve the paragraph index of the anchor, then grab the paragraph right
# after it (start of the content to delete) and the very last paragraph in
# the document (end of the content to delete).
$anchorIndex = $anchor.Paragraphs.Item(1).Index
$firstToDelete = $d.Paragraphs.Item($anchorIndex + 1)
$lastToDelete = $d.Paragraphs.Item($d.Paragraphs.Count)

# Delete the whole trailing range in one shot (Draft 2, Draft 3, and the
# trailing empty paragraph that followed them).
$deleteRange = $d.Range($firstToDelete.Range.Start, $lastToDelete.Range.End)
$deleteRange.Delete()
